# Apply crypto price/volume updates for Thu Jul  4 06:44:28 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        # Looks like a plain number to Excel auto-detection - force text storage
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue 'D2' '58.886.62'
Set-TextValue 'E2' '  -3.56%  '

Set-TextValue 'D3' '3.212.79'
Set-TextValue 'E3' '  -4.59%  '

Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.05%  '

Set-TextValue 'D5' '535.82'
Set-TextValue 'E5' '  -5.70%  '

Set-TextValue 'D6' '135.29'
Set-TextValue 'E6' '  -9.21%  '

Set-TextValue 'E7' '  -0.03%  '

Set-TextValue 'D8' '3.210.57'
Set-TextValue 'E8' '  -4.69%  '

Set-TextValue 'D9' '0.457'
Set-TextValue 'E9' '  -4.79%  '

Set-TextValue 'D10' '7.56'
Set-TextValue 'E10' '  -5.20%  '

Set-TextValue 'E11' '  -6.59%  '

Set-TextValue 'D12' '0.393'
Set-TextValue 'E12' '  -5.46%  '

Set-TextValue 'D13' '3.761.17'
Set-TextValue 'E13' '  -4.66%  '

Set-TextValue 'E14' '  -0.54%  '

Set-TextValue 'D15' '25.89'
Set-TextValue 'E15' '  -7.66%  '

Set-TextValue 'D16' '3.208.07'
Set-TextValue 'E16' '  -4.52%  '

Set-TextValue 'D17' '0.0000158'
Set-TextValue 'E17' '  -7.17%  '

Set-TextValue 'D18' '58.882.41'
Set-TextValue 'E18' '  -3.66%  '

Set-TextValue 'D19' '5.92'
Set-TextValue 'E19' '  -6.57%  '

Set-TextValue 'D20' '13.33'
Set-TextValue 'E20' '  -7.90%  '

Set-TextValue 'D21' '8.22'
Set-TextValue 'E21' '  -7.58%  '

Set-TextValue 'D22' '361.31'
Set-TextValue 'E22' '  -3.79%  '

Set-TextValue 'E23' '  -0.12%  '

Set-TextValue 'D24' '70.15'
Set-TextValue 'E24' '  -6.92%  '

Set-TextValue 'D25' '0.520'
Set-TextValue 'E25' '  -7.30%  '

Set-TextValue 'D26' '3.338.69'
Set-TextValue 'E26' '  -4.76%  '

Set-TextValue 'D27' '0.172'
Set-TextValue 'E27' '  -2.29%  '

Set-TextValue 'D28' '0.0₃0966'
Set-TextValue 'E28' '  -11.03%  '

Set-TextValue 'E29' '  +0.56%  '

Set-TextValue 'D30' '7.12'
Set-TextValue 'E30' '  -4.55%  '

Set-TextValue 'E31' '  +0.00%  '

Set-TextValue 'D32' '1.93'
Set-TextValue 'E32' '  -7.61%  '

Set-TextValue 'D33' '7.07'
Set-TextValue 'E33' '  -8.41%  '

Set-TextValue 'D34' '21.75'
Set-TextValue 'E34' '  -5.00%  '

Set-TextValue 'D35' '1.21'
Set-TextValue 'E35' '  -7.12%  '

Set-TextValue 'D36' '161.67'
Set-TextValue 'E36' '  -5.19%  '

Set-TextValue 'D37' '4.91'
Set-TextValue 'E37' '  -8.86%  '

Set-TextValue 'D38' '6.37'
Set-TextValue 'E38' '  -6.58%  '

Set-TextValue 'E39' '  -8.17%  '

Set-TextValue 'D40' '25.98'
Set-TextValue 'E40' '  -10.80%  '

Set-TextValue 'D41' '0.0707'
Set-TextValue 'E41' '  -6.42%  '

Set-TextValue 'D42' '3.239.67'
Set-TextValue 'E42' '  -4.87%  '

Set-TextValue 'D43' '40.95'
Set-TextValue 'E43' '  -3.46%  '

Set-TextValue 'D44' '0.715'
Set-TextValue 'E44' '  -6.15%  '

Set-TextValue 'B45' 'Filecoin'
Set-TextValue 'C45' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D45' '4.03'
Set-TextValue 'E45' '  -6.34%  '

Set-TextValue 'B46' 'ONDO'
Set-TextValue 'C46' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D46' '1.10'
Set-TextValue 'E46' '  -4.34%  '

Set-TextValue 'D47' '1.50'
Set-TextValue 'E47' '  -6.86%  '

Set-TextValue 'D48' '0.999'
Set-TextValue 'E48' '  -0.14%  '

Set-TextValue 'D49' '2.305.09'
Set-TextValue 'E49' '  -7.27%  '

Set-TextValue 'D50' '6.28'
Set-TextValue 'E50' '  -5.94%  '

Set-TextValue 'D51' '20.76'
Set-TextValue 'E51' '  -8.02%  '
